$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.852.06"
$ws.Range("E2").Value = "  -1.02%  "
$ws.Range("D3").Value = "1.735.74"
$ws.Range("E3").Value = "  -1.83%  "
$ws.Range("D4").Value = "'0.9997"
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").Value = "'229.11"
$ws.Range("E5").Value = "  -3.77%  "
$ws.Range("D6").Value = "'0.9996"
$ws.Range("E6").Value = "  -0.07%  "
$ws.Range("D7").Value = "'0.5260"
$ws.Range("E7").Value = "  +0.41%  "
$ws.Range("D8").Value = "'0.2756"
$ws.Range("E8").Value = "  -0.27%  "
$ws.Range("D9").Value = "'39.36"
$ws.Range("E9").Value = "  -3.03%  "
$ws.Range("D10").Value = "'0.06145"
$ws.Range("E10").Value = "  -1.23%  "
$ws.Range("D11").Value = "1.736.74"
$ws.Range("E11").Value = "  -1.78%  "
$ws.Range("D12").Value = "'0.07077"
$ws.Range("E12").Value = "  +0.77%  "
$ws.Range("E13").Value = "  -6.10%  "
$ws.Range("D14").Value = "'0.6408"
$ws.Range("E14").Value = "  -1.37%  "
$ws.Range("D15").Value = "'4.525"
$ws.Range("E15").Value = "  -0.06%  "
$ws.Range("D16").Value = "'76.81"
$ws.Range("E16").Value = "  -2.07%  "
$ws.Range("D17").Value = "'0.9997"
$ws.Range("E17").Value = "  -0.05%  "
$ws.Range("D18").Value = "'0.9994"
$ws.Range("E18").Value = "  -0.05%  "
$ws.Range("D19").Value = "25.830.24"
$ws.Range("E19").Value = "  -1.13%  "
$ws.Range("E20").Value = "  -2.04%  "
$ws.Range("D21").Value = "'0.000006656"
$ws.Range("E21").Value = "  -1.97%  "
$ws.Range("D22").Value = "1.959.64"
$ws.Range("E22").Value = "  -2.05%  "
$ws.Range("D23").Value = "'4.236"
$ws.Range("E23").Value = "  +3.56%  "
$ws.Range("D24").Value = "'8.789"
$ws.Range("E24").Value = "  +4.17%  "
$ws.Range("D25").Value = "'5.165"
$ws.Range("E25").Value = "  -0.87%  "
$ws.Range("E26").Value = "  +1.53%  "
$ws.Range("D27").Value = "'1.511"
$ws.Range("E27").Value = "  +1.49%  "
$ws.Range("D28").Value = "'15.14"
$ws.Range("E28").Value = "  -0.39%  "
$ws.Range("E29").Value = "  -4.55%  "
$ws.Range("D30").Value = "'102.22"
$ws.Range("E30").Value = "  -0.52%  "
$ws.Range("D31").Value = "'0.08328"
$ws.Range("E31").Value = "  -0.91%  "
$ws.Range("D32").Value = "'3.714"
$ws.Range("E32").Value = "  -0.59%  "
$ws.Range("D33").Value = "'3.535"
$ws.Range("E33").Value = "  +1.83%  "
$ws.Range("D34").Value = "'0.04482"
$ws.Range("E34").Value = "  +0.50%  "
$ws.Range("D35").Value = "'2.612"
$ws.Range("E35").Value = "  -1.64%  "
$ws.Range("E36").Value = "  -3.58%  "
$ws.Range("D37").Value = "'0.6190"
$ws.Range("E37").Value = "  +1.25%  "
$ws.Range("D38").Value = "'2.676"
$ws.Range("E38").Value = "  -3.39%  "
$ws.Range("D39").Value = "'0.01571"
$ws.Range("E39").Value = "  -1.02%  "
$ws.Range("B40").Value = "RenderToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D40").Value = "'1.907"
$ws.Range("E40").Value = "  -4.42%  "
$ws.Range("B41").Value = "PaxDollar"
$ws.Range("C41").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D41").Value = "'0.9992"
$ws.Range("E41").Value = "  -0.28%  "
$ws.Range("D42").Value = "'100.02"
$ws.Range("E42").Value = "  -2.96%  "
$ws.Range("D43").Value = "'0.3850"
$ws.Range("E43").Value = "  -0.92%  "
$ws.Range("D44").Value = "'5.031"
$ws.Range("E44").Value = "  +1.55%  "
$ws.Range("D45").Value = "'0.7248"
$ws.Range("E45").Value = "  -3.68%  "
$ws.Range("D46").Value = "'0.05324"
$ws.Range("E46").Value = "  -3.45%  "
$ws.Range("D47").Value = "'0.1121"
$ws.Range("E47").Value = "  -0.03%  "
$ws.Range("D48").Value = "'6.199"
$ws.Range("E48").Value = "  -3.84%  "
$ws.Range("D49").Value = "'53.33"
$ws.Range("E49").Value = "  +0.64%  "
$ws.Range("D50").Value = "'30.00"
$ws.Range("E50").Value = "  -1.24%  "
$ws.Range("D51").Value = "'7.630"
$ws.Range("E51").Value = "  +2.06%  "
